$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F, shifting the old "Generated Answer" column (F) to G.
$ws.Columns("F:F").Insert()

# New column header for the inserted column F.
$ws.Range("F1").Value = "Generated Result"

# Update D2 (Generated Query) text.
$ws.Range("D2").Value = 'SELECT city_name FROM city WHERE state_name = "arizona" ORDER BY population DESC LIMIT 1'

# Update Gold Result (E) values to the tuple-list representation, and copy
# the same values into the newly inserted Generated Result (F) column.
$ws.Range("E2").Value = "[('phoenix',)]"
$ws.Range("F2").Value = "[('phoenix',)]"

$ws.Range("E3").Value = "[('delaware',), ('allegheny',), ('hudson',)]"
$ws.Range("F3").Value = "[('delaware',), ('allegheny',), ('hudson',)]"

$ws.Range("E4").Value = "[(266807.0,)]"
$ws.Range("F4").Value = "[(266807.0,)]"

$ws.Range("E5").Value = "[(4113200,)]"
$ws.Range("F5").Value = "[(4113200,)]"
